$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# SingleInstance (row 4) now works -> flip status from 0 (red) to 1 (green)
$ws.Range("C4").Value = 1
$ws.Range("C3").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Clear the now-obsolete bug comment in D4
$ws.Range("D4").ClearContents()

# The HTML-generation-only caveat no longer applies separately; merge with
# the existing "no PDF" note used elsewhere in the sheet
$ws.Range("D20").Value = "Не формируется PDF (iTextSharp), dll-depend"

# Column D's "best fit" width shrinks now that the longest comment
# (the old SingleInstance / HTML-only notes) is gone
$ws.Columns.Item(4).ColumnWidth = 59.33

# Move the active selection to D5
$ws.Range("D5").Select()
